# LDLC price-history workbook: append a new scrape column (R) after the
# existing last column (Q). This mirrors how a new "run" of the scraper
# adds one more timestamped column to the sheet:
#   - R1 gets the new scrape timestamp, formatted like the other header cells
#   - R2:R100 get the same price that was just recorded in Q2:Q100
#     (no price changed since the previous run)
#   - R101:R204 stay blank (those products have no price yet), but the
#     cell is still materialized so every row spans through column R

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 100   # rows 2-100 hold product price data
$lastRow     = 204   # rows 101-204 exist but have no price data yet
$newCol      = 18    # column R
$srcCol      = 17    # column Q (the previous last column)

# --- Header: clone Q1's look (bold, centered, bordered) onto R1, then stamp
#     it with the new scrape date/time (kept as text, like its neighbours) ---
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("R1").Value = "2026-01-28 09:23:46"

# --- Data rows: copy this run's price (column Q) into the new column R ---
for ($r = 2; $r -le $lastDataRow; $r++) {
    $price = $ws.Cells.Item($r, $srcCol).Value2
    $ws.Cells.Item($r, $newCol).Value2 = $price
}

# --- Trailing rows: no price recorded for these products yet, but still
#     touch column R so the cell exists (kept blank) on every row ---
for ($r = ($lastDataRow + 1); $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newCol).Font.Bold = $false
}
